# Apply updated TPM-derived values to the LR-pairs sheet (Col9a3-Mag)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 7).Value = 1.345040666666667
$ws.Cells.Item(2, 8).Value = 4.035122
$ws.Cells.Item(2, 9).Value = 0.2185308326579933
$ws.Cells.Item(2, 10).Value = 0.2185308326579933
$ws.Cells.Item(2, 15).Value = 0.05546670559109387
$ws.Cells.Item(2, 16).Value = 0.05546670559109387
$ws.Cells.Item(2, 17).Value = 0.08245906311066667
$ws.Cells.Item(2, 18).Value = 0.742131567996
$ws.Cells.Item(2, 19).Value = 0.01212118535761752
$ws.Cells.Item(2, 20).Value = 0.01212118535761752

# Row 3
$ws.Cells.Item(3, 7).Value = 1.345040666666667
$ws.Cells.Item(3, 8).Value = 4.035122
$ws.Cells.Item(3, 9).Value = 0.2185308326579933
$ws.Cells.Item(3, 10).Value = 0.2185308326579933
$ws.Cells.Item(3, 13).Value = 0.5397903333333334
$ws.Cells.Item(3, 14).Value = 1.619371
$ws.Cells.Item(3, 15).Value = 0.488376202980433
$ws.Cells.Item(3, 16).Value = 0.4883762029804329
$ws.Cells.Item(3, 17).Value = 0.726039949806889
$ws.Cells.Item(3, 18).Value = 6.534359548262001
$ws.Cells.Item(3, 19).Value = 0.1067252582876632
$ws.Cells.Item(3, 20).Value = 0.1067252582876632

# Row 4
$ws.Cells.Item(4, 7).Value = 1.345040666666667
$ws.Cells.Item(4, 8).Value = 4.035122
$ws.Cells.Item(4, 9).Value = 0.2185308326579933
$ws.Cells.Item(4, 10).Value = 0.2185308326579933
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.005069
$ws.Cells.Item(4, 14).Value = 0.015207
$ws.Cells.Item(4, 15).Value = 0.004586186191257867
$ws.Cells.Item(4, 16).Value = 0.004586186191257867
$ws.Cells.Item(4, 17).Value = 0.006818011139333334
$ws.Cells.Item(4, 18).Value = 0.061362100254
$ws.Cells.Item(4, 19).Value = 0.001002223087100173
$ws.Cells.Item(4, 20).Value = 0.001002223087100173

# Row 5
$ws.Cells.Item(5, 7).Value = 1.345040666666667
$ws.Cells.Item(5, 8).Value = 4.035122
$ws.Cells.Item(5, 9).Value = 0.2185308326579933
$ws.Cells.Item(5, 10).Value = 0.2185308326579933
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.4991103333333333
$ws.Cells.Item(5, 14).Value = 1.497331
$ws.Cells.Item(5, 15).Value = 0.4515709052372154
$ws.Cells.Item(5, 16).Value = 0.4515709052372153
$ws.Cells.Item(5, 17).Value = 0.671323695486889
$ws.Cells.Item(5, 18).Value = 6.041913259382
$ws.Cells.Item(5, 19).Value = 0.09868216592561246
$ws.Cells.Item(5, 20).Value = 0.09868216592561245

# Row 6
$ws.Cells.Item(6, 9).Value = 0.1717131411100012
$ws.Cells.Item(6, 10).Value = 0.1717131411100013
$ws.Cells.Item(6, 15).Value = 0.05546670559109387
$ws.Cells.Item(6, 16).Value = 0.05546670559109387
$ws.Cells.Item(6, 19).Value = 0.009524362244070397
$ws.Cells.Item(6, 20).Value = 0.009524362244070397

# Row 7
$ws.Cells.Item(7, 9).Value = 0.1717131411100012
$ws.Cells.Item(7, 10).Value = 0.1717131411100013
$ws.Cells.Item(7, 13).Value = 0.5397903333333334
$ws.Cells.Item(7, 14).Value = 1.619371
$ws.Cells.Item(7, 15).Value = 0.488376202980433
$ws.Cells.Item(7, 16).Value = 0.4883762029804329
$ws.Cells.Item(7, 17).Value = 0.5704943272137778
$ws.Cells.Item(7, 18).Value = 5.134448944924001
$ws.Cells.Item(7, 19).Value = 0.08386061185714572
$ws.Cells.Item(7, 20).Value = 0.0838606118571457

# Row 8
$ws.Cells.Item(8, 9).Value = 0.1717131411100012
$ws.Cells.Item(8, 10).Value = 0.1717131411100013
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.005069
$ws.Cells.Item(8, 14).Value = 0.015207
$ws.Cells.Item(8, 15).Value = 0.004586186191257867
$ws.Cells.Item(8, 16).Value = 0.004586186191257867
$ws.Cells.Item(8, 17).Value = 0.005357331478666666
$ws.Cells.Item(8, 18).Value = 0.048215983308
$ws.Cells.Item(8, 19).Value = 0.0007875084366162013
$ws.Cells.Item(8, 20).Value = 0.0007875084366162014

# Row 9
$ws.Cells.Item(9, 9).Value = 0.1717131411100012
$ws.Cells.Item(9, 10).Value = 0.1717131411100013
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.4991103333333333
$ws.Cells.Item(9, 14).Value = 1.497331
$ws.Cells.Item(9, 15).Value = 0.4515709052372154
$ws.Cells.Item(9, 16).Value = 0.4515709052372153
$ws.Cells.Item(9, 17).Value = 0.5275003945737777
$ws.Cells.Item(9, 18).Value = 4.747503551164001
$ws.Cells.Item(9, 19).Value = 0.07754065857216896
$ws.Cells.Item(9, 20).Value = 0.07754065857216896

# Row 10
$ws.Cells.Item(10, 7).Value = 2.512500666666666
$ws.Cells.Item(10, 8).Value = 7.537502
$ws.Cells.Item(10, 9).Value = 0.4082098603762884
$ws.Cells.Item(10, 10).Value = 0.4082098603762884
$ws.Cells.Item(10, 15).Value = 0.05546670559109387
$ws.Cells.Item(10, 16).Value = 0.05546670559109387
$ws.Cells.Item(10, 17).Value = 0.1540313658706667
$ws.Cells.Item(10, 18).Value = 1.386282292836
$ws.Cells.Item(10, 19).Value = 0.02264205614487312
$ws.Cells.Item(10, 20).Value = 0.02264205614487312

# Row 11
$ws.Cells.Item(11, 7).Value = 2.512500666666666
$ws.Cells.Item(11, 8).Value = 7.537502
$ws.Cells.Item(11, 9).Value = 0.4082098603762884
$ws.Cells.Item(11, 10).Value = 0.4082098603762884
$ws.Cells.Item(11, 13).Value = 0.5397903333333334
$ws.Cells.Item(11, 14).Value = 1.619371
$ws.Cells.Item(11, 15).Value = 0.488376202980433
$ws.Cells.Item(11, 16).Value = 0.4883762029804329
$ws.Cells.Item(11, 17).Value = 1.356223572360222
$ws.Cells.Item(11, 18).Value = 12.206012151242
$ws.Cells.Item(11, 19).Value = 0.1993599816297444
$ws.Cells.Item(11, 20).Value = 0.1993599816297444

# Row 12
$ws.Cells.Item(12, 7).Value = 2.512500666666666
$ws.Cells.Item(12, 8).Value = 7.537502
$ws.Cells.Item(12, 9).Value = 0.4082098603762884
$ws.Cells.Item(12, 10).Value = 0.4082098603762884
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.005069
$ws.Cells.Item(12, 14).Value = 0.015207
$ws.Cells.Item(12, 15).Value = 0.004586186191257867
$ws.Cells.Item(12, 16).Value = 0.004586186191257867
$ws.Cells.Item(12, 17).Value = 0.01273586587933333
$ws.Cells.Item(12, 18).Value = 0.114622792914
$ws.Cells.Item(12, 19).Value = 0.001872126424793036
$ws.Cells.Item(12, 20).Value = 0.001872126424793036

# Row 13
$ws.Cells.Item(13, 7).Value = 2.512500666666666
$ws.Cells.Item(13, 8).Value = 7.537502
$ws.Cells.Item(13, 9).Value = 0.4082098603762884
$ws.Cells.Item(13, 10).Value = 0.4082098603762884
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.4991103333333333
$ws.Cells.Item(13, 14).Value = 1.497331
$ws.Cells.Item(13, 15).Value = 0.4515709052372154
$ws.Cells.Item(13, 16).Value = 0.4515709052372153
$ws.Cells.Item(13, 17).Value = 1.254015045240222
$ws.Cells.Item(13, 18).Value = 11.286135407162
$ws.Cells.Item(13, 19).Value = 0.1843356961768778
$ws.Cells.Item(13, 20).Value = 0.1843356961768778

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3216003333333333
$ws.Cells.Item(14, 8).Value = 0.9648009999999999
$ws.Cells.Item(14, 9).Value = 0.0522509024211076
$ws.Cells.Item(14, 10).Value = 0.0522509024211076
$ws.Cells.Item(14, 15).Value = 0.05546670559109387
$ws.Cells.Item(14, 16).Value = 0.05546670559109387
$ws.Cells.Item(14, 17).Value = 0.01971603003533333
$ws.Cells.Item(14, 18).Value = 0.177444270318
$ws.Cells.Item(14, 19).Value = 0.002898185421460549
$ws.Cells.Item(14, 20).Value = 0.002898185421460549

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3216003333333333
$ws.Cells.Item(15, 8).Value = 0.9648009999999999
$ws.Cells.Item(15, 9).Value = 0.0522509024211076
$ws.Cells.Item(15, 10).Value = 0.0522509024211076
$ws.Cells.Item(15, 13).Value = 0.5397903333333334
$ws.Cells.Item(15, 14).Value = 1.619371
$ws.Cells.Item(15, 15).Value = 0.488376202980433
$ws.Cells.Item(15, 16).Value = 0.4883762029804329
$ws.Cells.Item(15, 17).Value = 0.1735967511301111
$ws.Cells.Item(15, 18).Value = 1.562370760171
$ws.Cells.Item(15, 19).Value = 0.02551809732672165
$ws.Cells.Item(15, 20).Value = 0.02551809732672164

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3216003333333333
$ws.Cells.Item(16, 8).Value = 0.9648009999999999
$ws.Cells.Item(16, 9).Value = 0.0522509024211076
$ws.Cells.Item(16, 10).Value = 0.0522509024211076
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.005069
$ws.Cells.Item(16, 14).Value = 0.015207
$ws.Cells.Item(16, 15).Value = 0.004586186191257867
$ws.Cells.Item(16, 16).Value = 0.004586186191257867
$ws.Cells.Item(16, 17).Value = 0.001630192089666667
$ws.Cells.Item(16, 18).Value = 0.014671728807
$ws.Cells.Item(16, 19).Value = 0.0002396323671644459
$ws.Cells.Item(16, 20).Value = 0.0002396323671644459

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3216003333333333
$ws.Cells.Item(17, 8).Value = 0.9648009999999999
$ws.Cells.Item(17, 9).Value = 0.0522509024211076
$ws.Cells.Item(17, 10).Value = 0.0522509024211076
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.4991103333333333
$ws.Cells.Item(17, 14).Value = 1.497331
$ws.Cells.Item(17, 15).Value = 0.4515709052372154
$ws.Cells.Item(17, 16).Value = 0.4515709052372153
$ws.Cells.Item(17, 17).Value = 0.1605140495701111
$ws.Cells.Item(17, 18).Value = 1.444626446131
$ws.Cells.Item(17, 19).Value = 0.02359498730576097
$ws.Cells.Item(17, 20).Value = 0.02359498730576097

# Row 18
$ws.Cells.Item(18, 7).Value = 0.918901
$ws.Cells.Item(18, 8).Value = 2.756703
$ws.Cells.Item(18, 9).Value = 0.1492952634346094
$ws.Cells.Item(18, 10).Value = 0.1492952634346094
$ws.Cells.Item(18, 15).Value = 0.05546670559109387
$ws.Cells.Item(18, 16).Value = 0.05546670559109387
$ws.Cells.Item(18, 17).Value = 0.056334144706
$ws.Cells.Item(18, 18).Value = 0.5070073023539999
$ws.Cells.Item(18, 19).Value = 0.008280916423072283
$ws.Cells.Item(18, 20).Value = 0.008280916423072281

# Row 19
$ws.Cells.Item(19, 7).Value = 0.918901
$ws.Cells.Item(19, 8).Value = 2.756703
$ws.Cells.Item(19, 9).Value = 0.1492952634346094
$ws.Cells.Item(19, 10).Value = 0.1492952634346094
$ws.Cells.Item(19, 13).Value = 0.5397903333333334
$ws.Cells.Item(19, 14).Value = 1.619371
$ws.Cells.Item(19, 15).Value = 0.488376202980433
$ws.Cells.Item(19, 16).Value = 0.4883762029804329
$ws.Cells.Item(19, 17).Value = 0.4960138770903333
$ws.Cells.Item(19, 18).Value = 4.464124893813
$ws.Cells.Item(19, 19).Value = 0.07291225387915803
$ws.Cells.Item(19, 20).Value = 0.07291225387915802

# Row 20
$ws.Cells.Item(20, 7).Value = 0.918901
$ws.Cells.Item(20, 8).Value = 2.756703
$ws.Cells.Item(20, 9).Value = 0.1492952634346094
$ws.Cells.Item(20, 10).Value = 0.1492952634346094
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.3333333333333333
$ws.Cells.Item(20, 13).Value = 0.005069
$ws.Cells.Item(20, 14).Value = 0.015207
$ws.Cells.Item(20, 15).Value = 0.004586186191257867
$ws.Cells.Item(20, 16).Value = 0.004586186191257867
$ws.Cells.Item(20, 17).Value = 0.004657909169
$ws.Cells.Item(20, 18).Value = 0.041921182521
$ws.Cells.Item(20, 19).Value = 0.0006846958755840112
$ws.Cells.Item(20, 20).Value = 0.0006846958755840112

# Row 21
$ws.Cells.Item(21, 7).Value = 0.918901
$ws.Cells.Item(21, 8).Value = 2.756703
$ws.Cells.Item(21, 9).Value = 0.1492952634346094
$ws.Cells.Item(21, 10).Value = 0.1492952634346094
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.4991103333333333
$ws.Cells.Item(21, 14).Value = 1.497331
$ws.Cells.Item(21, 15).Value = 0.4515709052372154
$ws.Cells.Item(21, 16).Value = 0.4515709052372153
$ws.Cells.Item(21, 17).Value = 0.4586329844103333
$ws.Cells.Item(21, 18).Value = 4.127696859693
$ws.Cells.Item(21, 19).Value = 0.06741739725679512
$ws.Cells.Item(21, 20).Value = 0.0674173972567951

